$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so that numeric-looking
# strings like "1.002" are not auto-converted to numbers, matching the
# original inlineStr/text representation used throughout the sheet.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '28.139.03'
$ws.Range("E2").Value = '  -0.73%  '

# Row 3
$ws.Range("D3").Value = '1.826.25'
$ws.Range("E3").Value = '  +0.81%  '

# Row 4
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.24%  '

# Row 5
$ws.Range("D5").Value = '311.08'
$ws.Range("E5").Value = '  -0.73%  '

# Row 6
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.25%  '

# Row 7
$ws.Range("D7").Value = '0.4957'
$ws.Range("E7").Value = '  -3.84%  '

# Row 8
$ws.Range("D8").Value = '0.3926'
$ws.Range("E8").Value = '  -2.09%  '

# Row 9
$ws.Range("D9").Value = '0.09816'
$ws.Range("E9").Value = '  +24.64%  '

# Row 10
$ws.Range("D10").Value = '1.109'
$ws.Range("E10").Value = '  -0.67%  '

# Row 11
$ws.Range("D11").Value = '41.05'
$ws.Range("E11").Value = '  -0.02%  '

# Row 12
$ws.Range("D12").Value = '6.463'
$ws.Range("E12").Value = '  +1.27%  '

# Row 13
$ws.Range("D13").Value = '20.62'
$ws.Range("E13").Value = '  +0.66%  '

# Row 14
$ws.Range("D14").Value = '1.002'
$ws.Range("E14").Value = '  +0.28%  '

# Row 15
$ws.Range("D15").Value = '1.822.91'
$ws.Range("E15").Value = '  +1.29%  '

# Row 16
$ws.Range("D16").Value = '7.314'
$ws.Range("E16").Value = '  -0.70%  '

# Row 17
$ws.Range("D17").Value = '0.00001143'
$ws.Range("E17").Value = '  +5.45%  '

# Row 18
$ws.Range("D18").Value = '92.70'
$ws.Range("E18").Value = '  -0.25%  '

# Row 19
$ws.Range("D19").Value = '0.06664'
$ws.Range("E19").Value = '  +1.16%  '

# Row 20
$ws.Range("E20").Value = '  +0.16%  '

# Row 21
$ws.Range("D21").Value = '17.25'
$ws.Range("E21").Value = '  -0.73%  '

# Row 22
$ws.Range("D22").Value = '6.018'
$ws.Range("E22").Value = '  -0.43%  '

# Row 23
$ws.Range("D23").Value = '28.187.71'
$ws.Range("E23").Value = '  -0.74%  '

# Row 24
$ws.Range("D24").Value = '11.37'
$ws.Range("E24").Value = '  +1.37%  '

# Row 25
$ws.Range("D25").Value = '2.250'
$ws.Range("E25").Value = '  +0.95%  '

# Row 26
$ws.Range("D26").Value = '158.79'
$ws.Range("E26").Value = '  -1.21%  '

# Row 27
$ws.Range("D27").Value = '20.86'
$ws.Range("E27").Value = '  +1.07%  '

# Row 28
$ws.Range("D28").Value = '2.033.62'
$ws.Range("E28").Value = '  +0.83%  '

# Row 29
$ws.Range("D29").Value = '2.428'
$ws.Range("E29").Value = '  +0.39%  '

# Row 30
$ws.Range("D30").Value = '127.07'
$ws.Range("E30").Value = '  -1.28%  '

# Row 31
$ws.Range("D31").Value = '0.1055'
$ws.Range("E31").Value = '  -2.80%  '

# Row 32
$ws.Range("D32").Value = '1.040'
$ws.Range("E32").Value = '  -1.36%  '

# Row 33
$ws.Range("D33").Value = '5.608'
$ws.Range("E33").Value = '  +0.07%  '

# Row 34
$ws.Range("D34").Value = '3.616'
$ws.Range("E34").Value = '  -1.27%  '

# Row 35
$ws.Range("D35").Value = '0.06734'
$ws.Range("E35").Value = '  -6.51%  '

# Row 36
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.02355'
$ws.Range("E36").Value = '  +0.60%  '

# Row 37
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").Value = '8.989'
$ws.Range("E37").Value = '  -1.71%  '

# Row 38
$ws.Range("E38").Value = '  -0.86%  '

# Row 39
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '4.975'
$ws.Range("E39").Value = '  -2.04%  '

# Row 40
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '11.41'
$ws.Range("E40").Value = '  -2.16%  '

# Row 41
$ws.Range("D41").Value = '0.6229'
$ws.Range("E41").Value = '  -0.07%  '

# Row 42
$ws.Range("D42").Value = '1.181'
$ws.Range("E42").Value = '  +2.07%  '

# Row 43
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.26%  '

# Row 44
$ws.Range("D44").Value = '13.26'
$ws.Range("E44").Value = '  +0.23%  '

# Row 45
$ws.Range("D45").Value = '0.5936'
$ws.Range("E45").Value = '  -1.60%  '

# Row 46
$ws.Range("D46").Value = '3.711'
$ws.Range("E46").Value = '  -0.92%  '

# Row 47
$ws.Range("D47").Value = '1.278'
$ws.Range("E47").Value = '  -2.60%  '

# Row 48
$ws.Range("D48").Value = '124.28'
$ws.Range("E48").Value = '  -1.58%  '

# Row 49
$ws.Range("D49").Value = '1.952'
$ws.Range("E49").Value = '  +0.29%  '

# Row 50
$ws.Range("D50").Value = '1.185'
$ws.Range("E50").Value = '  -2.94%  '

# Row 51
$ws.Range("D51").Value = '0.06787'
$ws.Range("E51").Value = '  -1.02%  '
